# The document embeds the same two logo pictures twice (default +
# first-page header/footer). Their wp:docPr / pic:cNvPr "name"
# attributes were swapped between the pair of Pearson logos (footer)
# and bumped for the BTEC logo (header). InlineShape objects don't
# expose a settable Name in the Word OM for inline pictures, so the
# rename is done by editing the canonical package XML directly via
# Document.WordOpenXML (a round-trippable flat-OPC snapshot of the
# whole document, including headers/footers) and writing it back.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# --- Footer (first page) picture: Y:\...\PearsonLogo.png, docPr id="3" ---
# wp:docPr name="image2.png" -> name="image1.png"
$xml = $xml.Replace(
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image2.png"',
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"')

# --- Footer (default) picture: Y:\...\PearsonLogo.png, docPr id="2" ---
# wp:docPr name="image2.png" -> name="image1.png"
$xml = $xml.Replace(
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"',
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"')

# --- Both Pearson logos' pic:cNvPr (id="0") share identical text, so a
# single replace fixes both occurrences: name="image2.png" -> name="image1.png"
$xml = $xml.Replace(
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"',
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"')

# --- Header (first page) picture: BTec_Logo-Orange, docPr id="1" ---
# wp:docPr name="image1.jpg" -> name="image2.jpg"
$xml = $xml.Replace(
    'descr="BTec_Logo-Orange" id="1" name="image1.jpg"',
    'descr="BTec_Logo-Orange" id="1" name="image2.jpg"')

# --- Same picture's pic:cNvPr (id="0"): name="image1.jpg" -> name="image2.jpg" ---
$xml = $xml.Replace(
    'descr="BTec_Logo-Orange" id="0" name="image1.jpg"',
    'descr="BTec_Logo-Orange" id="0" name="image2.jpg"')

$d.WordOpenXML = $xml
